$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (not auto-coerced to a number) while
# preserving the cell's existing style. We do this by entering a formula
# that evaluates to the literal text, then collapsing it to a value via
# PasteSpecial(xlPasteValues) -- this keeps the shared-string cell type
# without pulling in a new number-format / style record (which a plain
# NumberFormat="@" trick, or a leading apostrophe, would do).
function Set-TextValue {
    param($range, $value)
    $range.Formula = '="' + $value + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# 1. Remove the discontinued line item (old row 9):
#    20139524 / MR.POTATO HDH PHOTO / R088 / 1 / 14 / RT
#    Rows 10-11 shift up to become rows 9-10.
$ws.Rows.Item(9).Delete()

# 2. Give the two brand-new rows (11 & 12) the same bordered cell style
#    used throughout the sheet, by copying formats from an existing row.
$ws.Range("A10:F10").Copy()
$ws.Range("A11:F12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3. New row 11: 20140198 / FIESTA TREASURE P/C / R088 / 4 / 2 / RT,(E-7H)
Set-TextValue $ws.Range("A11") "20140198"
$ws.Range("B11").Value2 = "FIESTA TREASURE P/C"
$ws.Range("C11").Value2 = "R088"
Set-TextValue $ws.Range("D11") "4"
Set-TextValue $ws.Range("E11") "2"
$ws.Range("F11").Value2 = "RT,(E-7H)"

# 4. New row 12: 20140668 / KIN PTCARD JKT48 SRS / R088 / 4 / 3 / RT
Set-TextValue $ws.Range("A12") "20140668"
$ws.Range("B12").Value2 = "KIN PTCARD JKT48 SRS"
$ws.Range("C12").Value2 = "R088"
Set-TextValue $ws.Range("D12") "4"
Set-TextValue $ws.Range("E12") "3"
$ws.Range("F12").Value2 = "RT"

# 5. Column F widens to fit the longer "RT,(E-7H)" value (xml width 11).
$ws.Columns.Item(6).ColumnWidth = 10.17
